$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.185315555808703
$ws.Range("C2").Value = 0.7782262909704238
$ws.Range("B3").Value = 9.447135311384551
$ws.Range("C3").Value = 1.95985241919762
$ws.Range("B4").Value = 15.64013449533723
$ws.Range("C4").Value = 2.911798748634491
$ws.Range("B5").Value = 15.71270320298209
$ws.Range("C5").Value = 3.78697830538002
$ws.Range("B6").Value = 15.97519774768478
$ws.Range("C6").Value = 4.795199348677913
$ws.Range("B7").Value = 16.04015616940114
$ws.Range("C7").Value = 5.761518494377109
$ws.Range("B8").Value = 19.75520056307449
$ws.Range("C8").Value = 6.748925997640637
$ws.Range("B9").Value = 23.84190300906377
$ws.Range("C9").Value = 7.668723920661449
$ws.Range("B10").Value = 23.87664518743228
$ws.Range("C10").Value = 8.665068007336934
$ws.Range("B11").Value = 29.25475426329459
$ws.Range("C11").Value = 9.394963830161716
$ws.Range("B12").Value = 30.06758303421456
$ws.Range("C12").Value = 10.17149312133259
$ws.Range("B13").Value = 31.02304977495041
$ws.Range("C13").Value = 11.05340859689823
$ws.Range("B14").Value = 31.13628047117407
$ws.Range("C14").Value = 12.00499148860196
$ws.Range("B15").Value = 33.37456298020494
$ws.Range("C15").Value = 12.90602332559623
$ws.Range("B16").Value = 38.96524419910486
$ws.Range("C16").Value = 13.83275445896211
$ws.Range("B17").Value = 39.74668882150542
$ws.Range("C17").Value = 14.64606333307143
$ws.Range("B18").Value = 40.01870186156726
$ws.Range("C18").Value = 15.53040120661702
$ws.Range("B19").Value = 42.76794407475322
$ws.Range("C19").Value = 16.43152165850896
$ws.Range("B20").Value = 43.77518365362869
$ws.Range("C20").Value = 17.23797235101107
$ws.Range("B21").Value = 45.14477737576481
$ws.Range("C21").Value = 18.06463496752176
$ws.Range("B22").Value = 45.21813687592143
$ws.Range("C22").Value = 18.97408342626731
$ws.Range("B23").Value = 45.28016142168637
$ws.Range("C23").Value = 19.88677315455734
$ws.Range("B24").Value = 45.47118925671833
$ws.Range("C24").Value = 21.06183242197405
$ws.Range("B25").Value = 45.54661958767856
$ws.Range("C25").Value = 21.95014255816888
$ws.Range("B26").Value = 49.13520885808486
$ws.Range("C26").Value = 22.82220535365735
$ws.Range("B27").Value = 50.92835473966905
$ws.Range("C27").Value = 23.75618307681371
$ws.Range("B28").Value = 51.26453784191647
$ws.Range("C28").Value = 24.69266577581023
$ws.Range("B29").Value = 53.03099358936512
$ws.Range("C29").Value = 25.61306973375901
$ws.Range("B30").Value = 56.6376471475108
$ws.Range("C30").Value = 26.53375034456346
$ws.Range("B31").Value = 68.9268166789534
$ws.Range("C31").Value = 27.40907657748353
$ws.Range("B32").Value = 70.925234315601
$ws.Range("C32").Value = 28.27044362247271
$ws.Range("B33").Value = 71.00669629517439
$ws.Range("C33").Value = 29.17203192327981
$ws.Range("B34").Value = 72.46640992221448
$ws.Range("C34").Value = 30.00774507921831
$ws.Range("B35").Value = 77.19908838582306
$ws.Range("C35").Value = 31.08034125291271
$ws.Range("B36").Value = 77.26955607505859
$ws.Range("C36").Value = 31.87670872747315
$ws.Range("B37").Value = 77.34336286177258
$ws.Range("C37").Value = 32.74389901949375
$ws.Range("B38").Value = 77.43956816434813
$ws.Range("C38").Value = 33.67043807320054
$ws.Range("B39").Value = 82.97519845047354
$ws.Range("C39").Value = 34.51558437015745
$ws.Range("B40").Value = 83.06741375479658
$ws.Range("C40").Value = 35.33346092965965
$ws.Range("B41").Value = 84.3530685895562
$ws.Range("C41").Value = 36.41652808240357
$ws.Range("B42").Value = 87.41360372121729
$ws.Range("C42").Value = 37.68930836904541
$ws.Range("B43").Value = 87.67977659055087
$ws.Range("C43").Value = 38.67798138596294
$ws.Range("B44").Value = 94.79699174739017
$ws.Range("C44").Value = 39.60189703367158
$ws.Range("B45").Value = 94.89318007233412
$ws.Range("C45").Value = 40.46169155816882
$ws.Range("B46").Value = 97.01794891088674
$ws.Range("C46").Value = 41.24420387961448
$ws.Range("B47").Value = 97.8336238251179
$ws.Range("C47").Value = 42.08466089909074
